# Weekly update: Fruta / hortaliza, semanal
# Insert two new observation rows (Black Amber "Especial" and "Primera"
# quality lots dated 2022-02-22) at rows 17-18 of the Ciruela sheet,
# pushing all the existing rows (old 17..41) down to rows 19..43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 17 downward by two, opening a blank 2-row gap at 17:18.
$ws.Range("A17:T18").Insert()

# --- New row 17 -----------------------------------------------------
$ws.Cells.Item(17, 1).Value = 7
$ws.Cells.Item(17, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(17, 3).Value = "Ñuble"
$ws.Cells.Item(17, 4).Value = 44614
$ws.Cells.Item(17, 5).Value = 16
$ws.Cells.Item(17, 6).Value = "Fruta"
$ws.Cells.Item(17, 7).Value = 100103
$ws.Cells.Item(17, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(17, 9).Value = 100103002
$ws.Cells.Item(17, 10).Value = "Ciruela"
$ws.Cells.Item(17, 11).Value = "Black Amber"
$ws.Cells.Item(17, 12).Value = "Especial"
$ws.Cells.Item(17, 13).Value = 50
$ws.Cells.Item(17, 14).Value = 13000
$ws.Cells.Item(17, 15).Value = 13000
$ws.Cells.Item(17, 16).Value = 13000
$ws.Cells.Item(17, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(17, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(17, 19).Value = 722
$ws.Cells.Item(17, 20).Value = 18

# --- New row 18 -----------------------------------------------------
$ws.Cells.Item(18, 1).Value = 7
$ws.Cells.Item(18, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(18, 3).Value = "Ñuble"
$ws.Cells.Item(18, 4).Value = 44614
$ws.Cells.Item(18, 5).Value = 16
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100103
$ws.Cells.Item(18, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(18, 9).Value = 100103002
$ws.Cells.Item(18, 10).Value = "Ciruela"
$ws.Cells.Item(18, 11).Value = "Black Amber"
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 13).Value = 120
$ws.Cells.Item(18, 14).Value = 11000
$ws.Cells.Item(18, 15).Value = 12000
$ws.Cells.Item(18, 16).Value = 11500
$ws.Cells.Item(18, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(18, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(18, 19).Value = 639
$ws.Cells.Item(18, 20).Value = 18
